$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections (the actual "fisher" data fix) ---
$ws.Range("C2").Value = 9
$ws.Range("D4").Value = 9
$ws.Range("D6").Value = 8

# --- Drop the now-unused helper column C in the blank filler rows ---
$ws.Range("C8:C20").Clear()

# --- Remove the two now-superfluous trailing blank rows ---
$ws.Rows("21:22").Delete()

# --- View state: zoom + selection ---
$ws.Range("D6").Select()
$excel.ActiveWindow.Zoom = 250
